$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = "Demo 1"

$ws.Range("G3").Value = 0.3
$ws.Range("H3").Value = "ha"
$ws.Range("K3").Value = "Demo 2"
